$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172504425048828
$ws.Range("B1").Value = 2.446961164474487
$ws.Range("C1").Value = 6.550660133361816
$ws.Range("D1").Value = 2.064251184463501
$ws.Range("E1").Value = 1.203955292701721
